$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 63,4
$arr[0,0] = 'word'
$arr[0,1] = 'meaning'
$arr[0,2] = 'example'
$arr[0,3] = 'example_english'
$arr[1,0] = 'adferden'
$arr[1,1] = 'the behavior'
$arr[1,2] = 'Adferden til barnet endret seg etter at han begynte på skolen.'
$arr[1,3] = 'The behavior of the child changed after he started school.'
$arr[2,0] = 'alvorlig'
$arr[2,1] = 'serious'
$arr[2,2] = 'Situasjonen er alvorlig, og vi må handle raskt.'
$arr[2,3] = 'The situation is serious, and we need to act quickly.'
$arr[3,0] = 'angriper'
$arr[3,1] = 'attacker'
$arr[3,2] = 'Løven angriper byttet sitt.'
$arr[3,3] = 'The lion attacks its prey.'
$arr[4,0] = 'avhengig'
$arr[4,1] = 'dependent'
$arr[4,2] = 'Han er avhengig av kaffe om morgenen.'
$arr[4,3] = 'He is dependent on coffee in the morning.'
$arr[5,0] = 'avslørt'
$arr[5,1] = 'revealed, exposed'
$arr[5,2] = 'Hans hemmelige plan ble avslørt.'
$arr[5,3] = 'His secret plan was revealed.'
$arr[6,0] = 'bedrager'
$arr[6,1] = 'fraudster'
$arr[6,2] = 'Politiet arresterte bedrageren som hadde stjålet tusenvis av kroner.'
$arr[6,3] = 'The police arrested the fraudster who had stolen thousands of kroner.'
$arr[7,0] = 'benyttet'
$arr[7,1] = 'used'
$arr[7,2] = 'Han benyttet anledningen til å takke alle involverte.'
$arr[7,3] = 'He used the opportunity to thank everyone involved.'
$arr[8,0] = 'betydning'
$arr[8,1] = 'meaning, significance'
$arr[8,2] = 'Ordet har en dyp betydning.'
$arr[8,3] = 'The word has a deep meaning.'
$arr[9,0] = 'bidrag'
$arr[9,1] = 'contribution'
$arr[9,2] = 'Hans bidrag til prosjektet var avgjørende for dets suksess.'
$arr[9,3] = 'His contribution to the project was crucial for its success.'
$arr[10,0] = 'bistand'
$arr[10,1] = 'assistance, aid'
$arr[10,2] = 'Vi tilbyr bistand til de som trenger det.'
$arr[10,3] = 'We offer assistance to those who need it.'
$arr[11,0] = 'blant'
$arr[11,1] = 'among'
$arr[11,2] = 'Han er populær blant sine venner.'
$arr[11,3] = 'He is popular among his friends.'
$arr[12,0] = 'byttet'
$arr[12,1] = 'exchanged, swapped, or the prey'
$arr[12,2] = 'Han byttet sin gamle bil mot en ny.'
$arr[12,3] = 'He exchanged his old car for a new one.'
$arr[13,0] = 'dekning'
$arr[13,1] = 'coverage'
$arr[13,2] = 'Det er dårlig dekning i dette området.'
$arr[13,3] = 'There is poor coverage in this area.'
$arr[14,0] = 'deretter'
$arr[14,1] = 'thereafter, then'
$arr[14,2] = 'Vi gikk hjem, og deretter så vi en film.'
$arr[14,3] = 'We went home, and then we watched a movie.'
$arr[15,0] = 'distraherte'
$arr[15,1] = 'distracted'
$arr[15,2] = 'Han distraherte meg mens jeg studerte.'
$arr[15,3] = 'He distracted me while I was studying.'
$arr[16,0] = 'dro'
$arr[16,1] = 'went, left'
$arr[16,2] = 'Han dro til butikken for å kjøpe melk.'
$arr[16,3] = 'He went to the store to buy milk.'
$arr[17,0] = 'egentlig'
$arr[17,1] = 'actually, really'
$arr[17,2] = 'Han er egentlig fra Norge, men bor i Sverige.'
$arr[17,3] = 'He is actually from Norway, but lives in Sweden.'
$arr[18,0] = 'enhet'
$arr[18,1] = 'unit'
$arr[18,2] = 'Hver enhet i bygningen har sin egen balkong.'
$arr[18,3] = 'Each unit in the building has its own balcony.'
$arr[19,0] = 'ensomme'
$arr[19,1] = 'lonely'
$arr[19,2] = 'Han følte seg ensomme i den store byen.'
$arr[19,3] = 'He felt lonely in the big city.'
$arr[20,0] = 'enten'
$arr[20,1] = 'either'
$arr[20,2] = 'Du må velge enten den røde eller den blå bilen.'
$arr[20,3] = 'You have to choose either the red car or the blue car.'
$arr[21,0] = 'erfaring'
$arr[21,1] = 'experience'
$arr[21,2] = 'Han har mye erfaring innen programmering.'
$arr[21,3] = 'He has a lot of experience in programming.'
$arr[22,0] = 'erstatte'
$arr[22,1] = 'replace'
$arr[22,2] = 'Vi må erstatte den gamle datamaskinen med en ny.'
$arr[22,3] = 'We need to replace the old computer with a new one.'
$arr[23,0] = 'ettermælet'
$arr[23,1] = 'legacy, reputation left behind'
$arr[23,2] = 'Hans ettermælet vil alltid bli husket for hans generøsitet.'
$arr[23,3] = 'His legacy will always be remembered for his generosity.'
$arr[24,0] = 'forbindelse'
$arr[24,1] = 'connection'
$arr[24,2] = 'Vi mistet forbindelsen under samtalen.'
$arr[24,3] = 'We lost the connection during the call.'
$arr[25,0] = 'forhandlingene'
$arr[25,1] = 'the negotiations'
$arr[25,2] = 'Forhandlingene mellom partene tok flere timer.'
$arr[25,3] = 'The negotiations between the parties took several hours.'
$arr[26,0] = 'forholde'
$arr[26,1] = 'relate, behave'
$arr[26,2] = 'Han vet ikke hvordan han skal forholde seg til situasjonen.'
$arr[26,3] = 'He doesn''t know how to behave in the situation.'
$arr[27,0] = 'fortjener'
$arr[27,1] = 'deserve'
$arr[27,2] = 'Du fortjener en belønning for alt ditt harde arbeid.'
$arr[27,3] = 'You deserve a reward for all your hard work.'
$arr[28,0] = 'hvilket som helst'
$arr[28,1] = 'whichever, any'
$arr[28,2] = 'Du kan velge hvilket som helst alternativ.'
$arr[28,3] = 'You can choose any option.'
$arr[29,0] = 'imponert'
$arr[29,1] = 'impressed'
$arr[29,2] = 'Jeg er veldig imponert over dine ferdigheter.'
$arr[29,3] = 'I am very impressed by your skills.'
$arr[30,0] = 'innsats'
$arr[30,1] = 'effort, contribution'
$arr[30,2] = 'Din innsats i prosjektet har vært uvurderlig.'
$arr[30,3] = 'Your contribution to the project has been invaluable.'
$arr[31,0] = 'kjedekollisjon'
$arr[31,1] = 'pile-up, multiple vehicle collision'
$arr[31,2] = 'Det var en stor kjedekollisjon på motorveien i går.'
$arr[31,3] = 'There was a big pile-up on the highway yesterday.'
$arr[32,0] = 'kjemper'
$arr[32,1] = 'fights, struggles, or giants (depending on context)'
$arr[32,2] = 'Han kjemper mot en vanskelig sykdom.'
$arr[32,3] = 'He is fighting against a difficult illness.'
$arr[33,0] = 'Klagar'
$arr[33,1] = 'Complains'
$arr[33,2] = 'Han klagar alltid på været.'
$arr[33,3] = 'He always complains about the weather.'
$arr[34,0] = 'krever'
$arr[34,1] = 'requires, demands'
$arr[34,2] = 'Denne oppgaven krever mye tid.'
$arr[34,3] = 'This task requires a lot of time.'
$arr[35,0] = 'likevel'
$arr[35,1] = 'nevertheless, nonetheless'
$arr[35,2] = 'Han var veldig trøtt, men han gikk på jobb likevel.'
$arr[35,3] = 'He was very tired, but he went to work nevertheless.'
$arr[36,0] = 'lyn'
$arr[36,1] = 'lightning'
$arr[36,2] = 'Vi så lynet slå ned i treet.'
$arr[36,3] = 'We saw the lightning strike the tree.'
$arr[37,0] = 'markedsføre'
$arr[37,1] = 'to market'
$arr[37,2] = 'Bedriften planlegger å markedsføre det nye produktet neste måned.'
$arr[37,3] = 'The company plans to market the new product next month.'
$arr[38,0] = 'minner'
$arr[38,1] = 'memories'
$arr[38,2] = 'Barndomsminnene mine er veldig kjære for meg.'
$arr[38,3] = 'My childhood memories are very dear to me.'
$arr[39,0] = 'naturligvis'
$arr[39,1] = 'naturally, of course'
$arr[39,2] = 'Naturligvis kan du låne bilen min.'
$arr[39,3] = 'Naturally, you can borrow my car.'
$arr[40,0] = 'ond'
$arr[40,1] = 'evil, wicked'
$arr[40,2] = 'Han hadde en ond plan.'
$arr[40,3] = 'He had an evil plan.'
$arr[41,0] = 'opplever'
$arr[41,1] = 'experiences'
$arr[41,2] = 'Hun opplever mye glede i sitt nye jobb.'
$arr[41,3] = 'She experiences a lot of joy in her new job.'
$arr[42,0] = 'overfladiske'
$arr[42,1] = 'superficial'
$arr[42,2] = 'Mange mennesker er bare interessert i overfladiske detaljer.'
$arr[42,3] = 'Many people are only interested in superficial details.'
$arr[43,0] = 'oversetter'
$arr[43,1] = 'translator'
$arr[43,2] = 'Han jobber som oversetter.'
$arr[43,3] = 'He works as a translator.'
$arr[44,0] = 'regjering'
$arr[44,1] = 'government'
$arr[44,2] = 'Regjeringen vil innføre nye skatteregler neste år.'
$arr[44,3] = 'The government will introduce new tax rules next year.'
$arr[45,0] = 'regning'
$arr[45,1] = 'bill, calculation'
$arr[45,2] = 'Kan jeg få regningen, vær så snill?'
$arr[45,3] = 'Can I have the bill, please?'
$arr[46,0] = 'rettferdighet'
$arr[46,1] = 'justice'
$arr[46,2] = 'Alle fortjener rettferdighet og like muligheter.'
$arr[46,3] = 'Everyone deserves justice and equal opportunities.'
$arr[47,0] = 'skrive inn'
$arr[47,1] = 'write in, enter'
$arr[47,2] = 'Kan du skrive inn navnet ditt her?'
$arr[47,3] = 'Can you write your name here?'
$arr[48,0] = 'spre'
$arr[48,1] = 'to spread'
$arr[48,2] = 'Han liker å spre glede blant vennene sine.'
$arr[48,3] = 'He likes to spread joy among his friends.'
$arr[49,0] = 'stjernerytter'
$arr[49,1] = 'star rider'
$arr[49,2] = 'Han drømte om å bli en stjernerytter og reise gjennom galaksene.'
$arr[49,3] = 'He dreamed of becoming a star rider and traveling through the galaxies.'
$arr[50,0] = 'svekker'
$arr[50,1] = 'weakens'
$arr[50,2] = 'Denne medisinen svekker immunforsvaret ditt.'
$arr[50,3] = 'This medicine weakens your immune system.'
$arr[51,0] = 'syn'
$arr[51,1] = 'sight, vision'
$arr[51,2] = 'Hans syn er ikke like godt som det en gang var.'
$arr[51,3] = 'His sight is not as good as it once was.'
$arr[52,0] = 'synes'
$arr[52,1] = 'to think, to seem, to appear'
$arr[52,2] = 'Jeg synes det er en god idé.'
$arr[52,3] = 'I think it is a good idea.'
$arr[53,0] = 'syns'
$arr[53,1] = 'to be seen, to seem, to appear'
$arr[53,2] = 'Det syns at han er glad.'
$arr[53,3] = 'It appears that he is happy.'
$arr[54,0] = 'trist'
$arr[54,1] = 'sad'
$arr[54,2] = 'Det er trist å se deg så lei deg.'
$arr[54,3] = 'It is sad to see you so upset.'
$arr[55,0] = 'turneringen'
$arr[55,1] = 'the tournament'
$arr[55,2] = 'Turneringen starter neste uke.'
$arr[55,3] = 'The tournament starts next week.'
$arr[56,0] = 'uansett'
$arr[56,1] = 'regardless'
$arr[56,2] = 'Vi går ut uansett været.'
$arr[56,3] = 'We are going out regardless of the weather.'
$arr[57,0] = 'Utfordrende'
$arr[57,1] = 'Challenging'
$arr[57,2] = 'Denne matteoppgaven er veldig utfordrende.'
$arr[57,3] = 'This math problem is very challenging.'
$arr[58,0] = 'uvurderlig'
$arr[58,1] = 'invaluable'
$arr[58,2] = 'Din støtte har vært uvurderlig for meg.'
$arr[58,3] = 'Your support has been invaluable to me.'
$arr[59,0] = 'vegne'
$arr[59,1] = 'on behalf of'
$arr[59,2] = 'Jeg snakker på vegne av hele teamet.'
$arr[59,3] = 'I speak on behalf of the entire team.'
$arr[60,0] = 'vekker'
$arr[60,1] = 'alarm clock'
$arr[60,2] = 'Jeg setter vekkeren på 7 om morgenen.'
$arr[60,3] = 'I set the alarm clock for 7 in the morning.'
$arr[61,0] = 'virkelig'
$arr[61,1] = 'real, really, actual'
$arr[61,2] = 'Er dette virkelig sant?'
$arr[61,3] = 'Is this really true?'
$arr[62,0] = 'årsaken'
$arr[62,1] = 'the reason'
$arr[62,2] = 'Årsaken til problemet er ikke kjent.'
$arr[62,3] = 'The reason for the problem is not known.'

$ws.Range("A1:D63").Value = $arr